# GradeTable: remove the final/as1/as2 assignment columns (F,G,H) and
# correct the course code from cn102 to cn101, per the "Add status and
# create_at field in GradeTable" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: clear the "final", "as1", "as2" labels in F1:H1 but keep
# their existing header style (s="1").
$ws.Range("F1").ClearContents()
$ws.Range("G1").ClearContents()
$ws.Range("H1").ClearContents()

# Data rows: drop the now-unused assignment score cells entirely.
$ws.Range("F2:H3").ClearContents()

# Fix the course code typo for both student rows.
$ws.Range("C2").Value = "cn101"
$ws.Range("C3").Value = "cn101"

# Match the author's final selection/view state.
[void]$ws.Range("F1:H3").Select()
